# Add a new "Yearly demand" worksheet at the end of the workbook, and
# populate it with the hourly profile data (mirrors the layout used by
# the other dispatch sheets: header row of hours 0-23 in B1:Y1, and a
# row label column A with day index 0,1,2 in rows 2-4).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Yearly demand"

# Match the page margins used throughout the rest of this workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 points)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Header row: hours 0..23 across columns B..Y (col index 2..25)
for ($c = 0; $c -le 23; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $c
}

# Row label column (A2:A4) = 0,1,2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2

$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

for ($c = 0; $c -le 23; $c++) {
    $ws.Cells.Item(2, $c + 2).Value = $row2[$c]
    $ws.Cells.Item(3, $c + 2).Value = $row3[$c]
    $ws.Cells.Item(4, $c + 2).Value = $row4[$c]
}

# Apply the same style used for headers/labels elsewhere in the workbook
# (bold, centered/top-aligned, thin box border) by copying formats from an
# existing header cell - this reuses the existing style entry instead of
# creating a near-duplicate one.
$styleSource = $wb.Worksheets.Item("Connected Households").Range("A2")
$styleSource.Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
